$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data keeps every Coin/Link/Price/Volume cell as plain text
# (inline strings), so force each touched cell to Text format before
# writing -- otherwise values such as "1.00" or "58.051.79" would be
# silently re-interpreted as numbers/dates and lose their formatting.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.053.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +3.05%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.065.68"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.34%  "

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +6.77%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.03"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +6.50%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.43%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.68"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.30%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.371"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +6.41%  "

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.51%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.589.22"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +3.27%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.22"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +8.83%  "

# Row 15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +17.59%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "58.051.79"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.23"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.06%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.066.26"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.19%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.15"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +6.68%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.68%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "339.74"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +4.58%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.03%  "

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.58%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.505"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.10"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.06%  "

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +6.68%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0980"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +9.73%  "

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.14%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.96"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.45%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +10.34%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.28%  "

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.33%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.21"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +5.39%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.77"

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.67"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.57%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.02"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +8.03%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.50%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.46"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +14.71%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0704"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.101.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.38%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.95"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.77%  "

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +10.41%  "

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.06%  "

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.667"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.82%  "

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.48"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.49%  "

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.338.41"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.18%  "

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.04%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.02"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +4.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0246"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +4.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.04"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +5.58%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.27"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +6.79%  "
